$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '26.741.08'
$ws.Range("E2").Value = '  -2.62%  '
$ws.Range("D3").Value = '1.565.94'
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '206.32'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.04%  '
$ws.Range("E6").Value = '  -2.35%  '
$ws.Range("E8").Value = '  -1.09%  '
$ws.Range("E9").Value = '  -0.87%  '
$ws.Range("E10").Value = '  -1.39%  '
$ws.Range("E11").Value = '  -0.35%  '
$ws.Range("D12").Value = '1.787.48'
$ws.Range("E12").Value = '  -0.12%  '
$ws.Range("D13").Value = '1.565.50'
$ws.Range("E13").Value = '  -0.27%  '
$ws.Range("E14").Value = '  -2.75%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.514'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -0.93%  '
$ws.Range("D16").Value = '26.783.35'
$ws.Range("E16").Value = '  -2.46%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.36'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -3.58%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '214.98'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +0.72%  '
$ws.Range("E19").Value = '  +1.51%  '
$ws.Range("D20").Value = '0.0₃0676'
$ws.Range("E20").Value = '  -1.94%  '
$ws.Range("E21").Value = '  +0.11%  '
$ws.Range("E23").Value = '  -2.73%  '
$ws.Range("E24").Value = '  -1.45%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '152.51'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.43%  '
$ws.Range("E26").Value = '  +0.49%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '14.92'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.38%  '
$ws.Range("E28").Value = '  +0.15%  '
$ws.Range("E29").Value = '  -1.48%  '
$ws.Range("E30").Value = '  -1.18%  '
$ws.Range("E31").Value = '  -3.62%  '
$ws.Range("E32").Value = '  -1.53%  '
$ws.Range("D33").Value = '1.389.41'
$ws.Range("E33").Value = '  +0.98%  '
$ws.Range("E34").Value = '  -1.33%  '
$ws.Range("E35").Value = '  +0.05%  '
$ws.Range("E36").Value = '  -0.93%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.929'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -2.67%  '
$ws.Range("E38").Value = '  -2.86%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.527'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -1.17%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.818'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -0.58%  '
$ws.Range("E41").Value = '  +0.11%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.991'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +1.41%  '
$ws.Range("E43").Value = '  -0.88%  '
$ws.Range("E44").Value = '  +1.36%  '
$ws.Range("E45").Value = '  +0.77%  '
$ws.Range("E46").Value = '  -1.53%  '
$ws.Range("D47").Value = '1.701.51'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '85.68'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.29%  '
$ws.Range("D49").Value = '0.0₇0981'
$ws.Range("E49").Value = '  -2.86%  '
$ws.Range("E50").Value = '  -0.87%  '
$ws.Range("E51").Value = '  -0.82%  '
